# chỉnh sửa một số chi tiết 13:34 25/03/21
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Màn hình chính" (sheet1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Màn hình chính")

# Insert a new blank row before the old row 6 ("Nghiệp vụ" row), shifting
# everything below down by one. The inserted row picks up an inherited
# style on the touched cell, so clear it to leave a truly empty row.
$ws1.Rows("6").Insert()
$ws1.Cells.Item(6, 2).Clear()

# New header row (STT / Testcase / Steps / Data Test / Kết quả mong muốn /
# Kết quả thực tế) - now spans 6 columns instead of 5.
$ws1.Range("A1").Value = "STT"
$ws1.Range("B1").Value = "Testcase"
$ws1.Range("C1").Value = "Steps"
$ws1.Range("D1").Value = "Data Test"
$ws1.Range("E1").Value = "Kết quả mong muốn"
$ws1.Range("F1").Value = "Kết quả thực tế"

# Row 8 (old row 7, shifted by the insert above) gets new wording and a
# taller custom row height.
$ws1.Range("E8").Value = "Mở form thêm mới, data trắng"
$ws1.Rows("8").RowHeight = 27.75

# New column F width.
$ws1.Columns("F").ColumnWidth = 36.15

# Restore the view's selection (best effort).
$ws1.Range("G8").Select()

# ---------------------------------------------------------------------
# Sheet "Màn hình thêm mới" (sheet2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Màn hình thêm mới")
$ws2.Range("D15").Select()

# ---------------------------------------------------------------------
# Sheet "Màn hình sửa" (sheet3)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Màn hình sửa")
